$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the typo in C18: "Clip 1-Bmode + Color Doppler" -> "Clip 1 B-mode + Color Doppler"
$ws.Range("C18").Value = "Clip 1 B-mode + Color Doppler"

# 2) Insert a new row at 19 for a second "Portal vein thrombosis" clip (shifts rows 19-26 -> 20-27)
$ws.Range("A19").EntireRow.Insert()

$ws.Range("A19").Value = "Liver vasculature"
$ws.Range("B19").Value = "Portal vein thrombosis"
$ws.Range("C19").Value = "Clip 2 B-mode + Color"
$ws.Range("D19").Value = "https://youtu.be/A3XUYC74J0o"
$ws.Range("D19").Style = "Collegamento ipertestuale"

# 3) The hyperlink objects on column D do not auto-shift with the row insert in this
#    environment, so rebuild them (in the original order) pointing at their
#    post-insert cell locations.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("D3"), "https://youtu.be/zxTC0YBY2RY")
$ws.Hyperlinks.Add($ws.Range("D21"), "https://youtu.be/xBfd04F4Ni8")
$ws.Hyperlinks.Add($ws.Range("D9"), "https://youtu.be/91M82AIMyu0")
$ws.Hyperlinks.Add($ws.Range("D25"), "https://youtu.be/qushjTAy6XQ")
$ws.Hyperlinks.Add($ws.Range("D23"), "https://youtu.be/pc-vbxSRTbs")
$ws.Hyperlinks.Add($ws.Range("D16"), "https://youtu.be/DjI1kEnzfSQ")
$ws.Hyperlinks.Add($ws.Range("D22"), "https://youtu.be/JvwODCASLYQ")
$ws.Hyperlinks.Add($ws.Range("D17"), "https://youtu.be/U3ydTsRwxok")
$ws.Hyperlinks.Add($ws.Range("D11"), "https://youtu.be/15o_Km86IzM")
$ws.Hyperlinks.Add($ws.Range("D26"), "https://youtu.be/_FckFwJwynI")
$ws.Hyperlinks.Add($ws.Range("D24"), "https://youtu.be/Axbee4vjNtU")
$ws.Hyperlinks.Add($ws.Range("D13"), "https://youtu.be/RhSUFLTmTl4")
$ws.Hyperlinks.Add($ws.Range("D6"), "https://youtu.be/2kRZcpi70Aw")
$ws.Hyperlinks.Add($ws.Range("D27"), "https://youtu.be/z_oaRVxRz5s")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://youtu.be/K2Wbg7BgXy4")

# Restore the plain "hyperlink" cell style (Add() above tends to create its own
# style variant) on every column-D cell that carries that formatting.
$ws.Range("D3").Style = "Collegamento ipertestuale"
$ws.Range("D4").Style = "Collegamento ipertestuale"
$ws.Range("D5").Style = "Collegamento ipertestuale"
$ws.Range("D6").Style = "Collegamento ipertestuale"
$ws.Range("D7").Style = "Collegamento ipertestuale"
$ws.Range("D8").Style = "Collegamento ipertestuale"
$ws.Range("D9").Style = "Collegamento ipertestuale"
$ws.Range("D10").Style = "Collegamento ipertestuale"
$ws.Range("D11").Style = "Collegamento ipertestuale"
$ws.Range("D12").Style = "Collegamento ipertestuale"
$ws.Range("D13").Style = "Collegamento ipertestuale"
$ws.Range("D14").Style = "Collegamento ipertestuale"
$ws.Range("D15").Style = "Collegamento ipertestuale"
$ws.Range("D16").Style = "Collegamento ipertestuale"
$ws.Range("D17").Style = "Collegamento ipertestuale"
$ws.Range("D18").Style = "Collegamento ipertestuale"
$ws.Range("D19").Style = "Collegamento ipertestuale"
$ws.Range("D21").Style = "Collegamento ipertestuale"
$ws.Range("D22").Style = "Collegamento ipertestuale"
$ws.Range("D23").Style = "Collegamento ipertestuale"
$ws.Range("D24").Style = "Collegamento ipertestuale"
$ws.Range("D25").Style = "Collegamento ipertestuale"
$ws.Range("D26").Style = "Collegamento ipertestuale"
$ws.Range("D27").Style = "Collegamento ipertestuale"

# Update selection to mirror the author's final cursor position
$ws.Range("D19").Select()
